$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new "Tipo" column is spliced in at G (only for the lower
# "crearSerVivo" table, rows 20-37). The pre-existing G formula
# (CONCAT(...)) moves to H, and the pre-existing H formula (trailing comma
# wrapper) moves to I - written directly cell-by-cell so the big matrix
# table above (rows 2-18, which already uses columns through Q) is left
# untouched.
# ---------------------------------------------------------------------------

# New header for the inserted "Tipo" column
$ws.Range("F20").HorizontalAlignment = -4108
$ws.Range("G20").Value = "Tipo"
$ws.Range("G20").HorizontalAlignment = -4108

# Animal -> Tipo (type) mapping for rows 21-36
$tipos = @{
  21 = "Carnivoro"   # LOBO
  22 = "Carnivoro"   # BOA
  23 = "Carnivoro"   # ZORRO
  24 = "Carnivoro"   # OSO
  25 = "Carnivoro"   # AGUILA
  26 = "Herbivoro"   # CABALLO
  27 = "Herbivoro"   # CIERVO
  28 = "Herbivoro"   # CONEJO
  29 = "Omnivoro"    # RATON
  30 = "Herbivoro"   # CABRA
  31 = "Herbivoro"   # OVEJA
  32 = "Omnivoro"    # JABALI
  33 = "Herbivoro"   # BUFALO
  34 = "Omnivoro"    # PATO
  35 = "Herbivoro"   # ORUGA
  36 = "Planta"      # PLANTA
}

foreach ($r in 21..36) {
  $tipo = $tipos[$r]

  # G: the new "Tipo" value for this row's species
  $ws.Range("G$r").Value = $tipo

  # H: rebuild the crearSerVivo(...) CONCAT formula (previously in G) so it
  # also appends the new Tipo argument (quoted) at the end of the call.
  $ws.Range("H$r").Formula = '=CONCAT(A' + $r + ',"(",CHAR(34),F' + $r + ',CHAR(34),",",B' + $r + ',",",C' + $r + ',",",D' + $r + ',",",E' + $r + ',",",CHAR(34),G' + $r + ',CHAR(34),")")'

  # I: trailing-comma wrapper around H (previously in H, around G), same
  # shape as before, just shifted one column right.
  $ws.Range("I$r").Formula = '=CONCAT(H' + $r + ',",")'
}

# ---------------------------------------------------------------------------
# Column widths: give the new G a width matching column A, and size the new
# H (formula) column a bit wider to fit the longer generated text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 34
$ws.Columns.Item(8).ColumnWidth = 41

# ---------------------------------------------------------------------------
# A handful of rows grow taller once the sheet is re-rendered with the new
# column in place (matches the committed workbook).
# ---------------------------------------------------------------------------
$tallRows = @(26, 28, 29, 33, 35)
foreach ($r in $tallRows) {
  $ws.Rows.Item($r).RowHeight = 24.05
}

# ---------------------------------------------------------------------------
# AutoFilter / _FilterDatabase range grows by one row (A2:Q17 -> A2:Q18).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A2:Q18").AutoFilter()
foreach ($n in $wb.Names) {
  if ($n.Name() -like "*_FilterDatabase*") {
    $n.RefersTo = "=Sheet1!`$A`$2:`$Q`$18"
  }
}

# ---------------------------------------------------------------------------
# Selection / scroll position, matching the new activity on row 34.
# ---------------------------------------------------------------------------
$ws.Range("H34").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
